$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 2 "Bitcoin"
Set-TextCell 2 3 "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell 2 4 "25.816.96"
Set-TextCell 2 5 "  +0.49%  "

Set-TextCell 3 2 "Ethereum"
Set-TextCell 3 3 "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell 3 4 "1.758.76"
Set-TextCell 3 5 "  -0.36%  "

Set-TextCell 4 2 "TetherUSD"
Set-TextCell 4 3 "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell 4 4 "1.000"
Set-TextCell 4 5 "  -0.53%  "

Set-TextCell 5 2 "BNB"
Set-TextCell 5 3 "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell 5 4 "236.96"
Set-TextCell 5 5 "  -0.16%  "

Set-TextCell 6 2 "USDC"
Set-TextCell 6 3 "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell 6 4 "0.9999"
Set-TextCell 6 5 "  -0.60%  "

Set-TextCell 7 2 "XRP"
Set-TextCell 7 3 "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell 7 4 "0.5062"
Set-TextCell 7 5 "  +3.03%  "

Set-TextCell 8 2 "OKB"
Set-TextCell 8 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 8 4 "41.28"
Set-TextCell 8 5 "  -2.29%  "

Set-TextCell 9 2 "Cardano"
Set-TextCell 9 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell 9 4 "0.2666"
Set-TextCell 9 5 "  +7.98%  "

Set-TextCell 10 2 "Dogecoin"
Set-TextCell 10 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 10 4 "0.06209"
Set-TextCell 10 5 "  +2.66%  "

Set-TextCell 11 2 "WrappedEther"
Set-TextCell 11 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 11 4 "1.755.26"
Set-TextCell 11 5 "  -0.58%  "

Set-TextCell 12 2 "Solana"
Set-TextCell 12 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 12 4 "15.68"
Set-TextCell 12 5 "  +9.70%  "

Set-TextCell 13 2 "TRON"
Set-TextCell 13 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 13 4 "0.06935"
Set-TextCell 13 5 "  +4.97%  "

Set-TextCell 14 2 "Polygon"
Set-TextCell 14 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 14 4 "0.6049"
Set-TextCell 14 5 "  +0.23%  "

Set-TextCell 15 2 "Polkadot"
Set-TextCell 15 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 15 4 "4.471"
Set-TextCell 15 5 "  +3.55%  "

Set-TextCell 16 2 "Litecoin"
Set-TextCell 16 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 16 4 "77.60"
Set-TextCell 16 5 "  -0.26%  "

Set-TextCell 17 2 "BinanceUSD"
Set-TextCell 17 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 17 4 "0.9999"
Set-TextCell 17 5 "  -0.61%  "

Set-TextCell 18 2 "Dai"
Set-TextCell 18 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 18 4 "0.9999"
Set-TextCell 18 5 "  -0.55%  "

Set-TextCell 19 2 "WrappedBTC"
Set-TextCell 19 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 19 4 "25.851.85"
Set-TextCell 19 5 "  +0.56%  "

Set-TextCell 20 2 "ShibaInu"
Set-TextCell 20 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 20 4 "0.000006841"
Set-TextCell 20 5 "  +8.85%  "

Set-TextCell 21 2 "Avalanche"
Set-TextCell 21 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 21 4 "11.60"
Set-TextCell 21 5 "  +5.41%  "

Set-TextCell 22 2 "WrappedliquidstakedEther2.0"
Set-TextCell 22 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 22 4 "1.977.10"
Set-TextCell 22 5 "  -1.05%  "

Set-TextCell 23 2 "Uniswap"
Set-TextCell 23 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 23 4 "4.069"
Set-TextCell 23 5 "  +5.57%  "

Set-TextCell 24 2 "Cosmos"
Set-TextCell 24 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 24 4 "8.185"
Set-TextCell 24 5 "  +2.15%  "

Set-TextCell 25 2 "Chainlink"
Set-TextCell 25 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 25 4 "5.190"
Set-TextCell 25 5 "  +1.45%  "

Set-TextCell 26 2 "Monero"
Set-TextCell 26 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 26 4 "138.15"
Set-TextCell 26 5 "  +4.41%  "

Set-TextCell 27 2 "Toncoin"
Set-TextCell 27 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 27 4 "1.455"
Set-TextCell 27 5 "  +5.69%  "

Set-TextCell 28 2 "LidoDAOToken"
Set-TextCell 28 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 28 4 "1.819"
Set-TextCell 28 5 "  -2.62%  "

Set-TextCell 29 2 "EthereumClassic"
Set-TextCell 29 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 29 4 "15.02"
Set-TextCell 29 5 "  +4.53%  "

Set-TextCell 30 2 "BitcoinCash"
Set-TextCell 30 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 30 4 "102.62"
Set-TextCell 30 5 "  +3.37%  "

Set-TextCell 31 2 "Stellar"
Set-TextCell 31 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 31 4 "0.08208"
Set-TextCell 31 5 "  +0.02%  "

Set-TextCell 32 2 "InternetComputer(DFINITY)"
Set-TextCell 32 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 32 4 "3.690"
Set-TextCell 32 5 "  +3.33%  "

Set-TextCell 33 2 "Filecoin"
Set-TextCell 33 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "3.392"
Set-TextCell 33 5 "  +7.38%  "

Set-TextCell 34 2 "Hedera"
Set-TextCell 34 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 34 4 "0.04383"
Set-TextCell 34 5 "  +2.31%  "

Set-TextCell 35 2 "Frax"
Set-TextCell 35 3 "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell 35 4 "0.9989"
Set-TextCell 35 5 "  -0.62%  "

Set-TextCell 36 2 "HuobiToken"
Set-TextCell 36 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 36 4 "2.650"
Set-TextCell 36 5 "  +1.09%  "

Set-TextCell 37 2 "ARBITRUM"
Set-TextCell 37 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 37 4 "0.9993"
Set-TextCell 37 5 "  -2.18%  "

Set-TextCell 38 2 "ImmutableX"
Set-TextCell 38 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 38 4 "0.6045"
Set-TextCell 38 5 "  -1.40%  "

Set-TextCell 39 2 "MXToken"
Set-TextCell 39 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 39 4 "2.732"
Set-TextCell 39 5 "  +2.71%  "

Set-TextCell 40 2 "VeChain"
Set-TextCell 40 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 40 4 "0.01547"
Set-TextCell 40 5 "  +7.85%  "

Set-TextCell 41 2 "RenderToken"
Set-TextCell 41 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 41 4 "1.938"
Set-TextCell 41 5 "  -7.98%  "

Set-TextCell 42 2 "PaxDollar"
Set-TextCell 42 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 42 4 "0.9999"
Set-TextCell 42 5 "  -0.60%  "

Set-TextCell 43 2 "Quant"
Set-TextCell 43 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell 43 4 "103.28"
Set-TextCell 43 5 "  +1.87%  "

Set-TextCell 44 2 "TheSandbox"
Set-TextCell 44 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 44 4 "0.3820"
Set-TextCell 44 5 "  +0.08%  "

Set-TextCell 45 2 "TrustWalletToken"
Set-TextCell 45 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 45 4 "0.7382"
Set-TextCell 45 5 "  -5.93%  "

Set-TextCell 46 2 "FraxShare"
Set-TextCell 46 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 46 4 "4.912"
Set-TextCell 46 5 "  -4.92%  "

Set-TextCell 47 2 "Cronos"
Set-TextCell 47 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell 47 4 "0.05491"
Set-TextCell 47 5 "  +6.03%  "

Set-TextCell 48 2 "Algorand"
Set-TextCell 48 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 48 4 "0.1086"
Set-TextCell 48 5 "  +6.76%  "

Set-TextCell 49 2 "Aptos"
Set-TextCell 49 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 49 4 "5.942"
Set-TextCell 49 5 "  -2.58%  "

Set-TextCell 50 2 "Elrond"
Set-TextCell 50 3 "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell 50 4 "29.80"
Set-TextCell 50 5 "  +2.89%  "

Set-TextCell 51 2 "EnergySwap"
Set-TextCell 51 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 51 4 "7.606"
Set-TextCell 51 5 "  +2.89%  "
